$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Text fix: report 0804 - note #6 is now marked as fixed ---
$cell = $ws.Range("I11")
$text = $cell.Text
$fixedText = $text.Replace(
    "6. Missing a condition to handle NOFUNDS, when the client does not have enough money to buy a ticket.",
    "6. Fixed: Missing a condition to handle NOFUNDS, when the client does not have enough money to buy a ticket."
)
$cell.Value = $fixedText

# --- Formatting cleanup: these cells drop their leftover "alignment" style back to Normal ---
$normalCells = @("D2","C3","D3","D4","D5","D6","D7","I7","D8","C9","D9","I9","C10","D10","I10","C11","D11")
foreach ($addr in $normalCells) {
    $ws.Range($addr).Style = "Normal"
}

# --- Formatting: I8 and I11 pick up the same Times New Roman formatting already used on C6 ---
$ws.Range("C6").Copy() | Out-Null
$ws.Range("I8").PasteSpecial(-4122) | Out-Null
$ws.Range("I11").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Restore the cursor/selection position as last left by the editor ---
$ws.Range("I13").Select()
